# Updates generated at 456a3b4
# Applies numeric refreshes across 展览/演出/本地生活/全部类型 sheets and
# appends/inserts the new "广州·火影only" event row.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

function Set-IndexCell {
    param($ws, [string]$addr, [string]$styleSourceAddr, $value)
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 47
$wsExpo.Range("F5").Value = 1023
$wsExpo.Range("G5").Value = 54
$wsExpo.Range("F7").Value = 609
$wsExpo.Range("G8").Value = 49.5
$wsExpo.Range("F9").Value = 1484
$wsExpo.Range("F12").Value = 3033
$wsExpo.Range("F13").Value = 506
$wsExpo.Range("F14").Value = 1687
$wsExpo.Range("F18").Value = 1421
$wsExpo.Range("F21").Value = 1154
$wsExpo.Range("F22").Value = 24
$wsExpo.Range("F24").Value = 32
$wsExpo.Range("F25").Value = 3584
$wsExpo.Range("F27").Value = 565
$wsExpo.Range("F28").Value = 1582

# New row 29 (appended at the end of the sheet)
Set-IndexCell $wsExpo "A29" "A28" 28
Set-TextCell  $wsExpo "B29" "2024-07-14"
Set-TextCell  $wsExpo "C29" "广州·火影only"
Set-TextCell  $wsExpo "D29" "人和镇蚌湖清河大街168号 人和园"
Set-TextCell  $wsExpo "E29" "2024.07.14 09:30-07.14 17:30"
$wsExpo.Range("F29").Value = 0
$wsExpo.Range("G29").Value = 78
Set-TextCell  $wsExpo "H29" "https://show.bilibili.com/platform/detail.html?id=84815"
Set-TextCell  $wsExpo "I29" "//i2.hdslb.com/bfs/openplatform/202404/nXYU1nYl1714035073709.png"

# ---------------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("F7").Value = 6
$wsShow.Range("F9").Value = 25

# ---------------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")

$wsLocal.Range("F2").Value = 801

# ---------------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 47
$wsAll.Range("F3").Value = 801
$wsAll.Range("F12").Value = 6
$wsAll.Range("F14").Value = 25
$wsAll.Range("F16").Value = 1023
$wsAll.Range("G16").Value = 54
$wsAll.Range("F18").Value = 609
$wsAll.Range("G19").Value = 49.5
$wsAll.Range("F20").Value = 1484
$wsAll.Range("F23").Value = 3033
$wsAll.Range("F24").Value = 506
$wsAll.Range("F25").Value = 1687
$wsAll.Range("F29").Value = 1421
$wsAll.Range("F34").Value = 1154
$wsAll.Range("F35").Value = 24
$wsAll.Range("F37").Value = 32
$wsAll.Range("F38").Value = 3584
$wsAll.Range("F40").Value = 565
$wsAll.Range("F41").Value = 1582

# Insert a new row before the existing row 44 ("孟京辉..." shifts to row 45)
# and populate the freed row 44 with the "广州·火影only" event.
$wsAll.Rows.Item(44).Insert()

Set-IndexCell $wsAll "A44" "A45" 43
$wsAll.Range("A45").Value = 44

Set-TextCell $wsAll "B44" "2024-07-14"
Set-TextCell $wsAll "C44" "广州·火影only"
Set-TextCell $wsAll "D44" "人和镇蚌湖清河大街168号 人和园"
Set-TextCell $wsAll "E44" "2024.07.14 09:30-07.14 17:30"
$wsAll.Range("F44").Value = 0
$wsAll.Range("G44").Value = 78
Set-TextCell $wsAll "H44" "https://show.bilibili.com/platform/detail.html?id=84815"
Set-TextCell $wsAll "I44" "//i2.hdslb.com/bfs/openplatform/202404/nXYU1nYl1714035073709.png"

"edits applied"
